$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '25.816.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = "'" + '1.738.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = "'" + '225.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.13%  '
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = "'" + '0.5148'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.97%  '
$ws.Range('D8').Value = "'" + '0.2716'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.87%  '
$ws.Range('D9').Value = "'" + '38.86'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.58%  '
$ws.Range('D10').Value = "'" + '0.06083'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('D11').Value = "'" + '1.740.18'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = "'" + '0.06993'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').Value = "'" + '15.25'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').Value = "'" + '0.6310'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.14%  '
$ws.Range('D15').Value = "'" + '4.486'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').Value = "'" + '76.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = "'" + '25.839.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = "'" + '11.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('D21').Value = "'" + '0.000006573'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').Value = "'" + '1.960.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = "'" + '4.035'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').Value = "'" + '8.410'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.46%  '
$ws.Range('D25').Value = "'" + '5.089'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').Value = "'" + '135.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').Value = "'" + '1.504'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.79%  '
$ws.Range('D28').Value = "'" + '1.816'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('D29').Value = "'" + '14.96'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').Value = "'" + '102.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('D31').Value = "'" + '0.08306'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').Value = "'" + '3.610'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').Value = "'" + '3.360'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.31%  '
$ws.Range('D34').Value = "'" + '0.04404'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('D35').Value = "'" + '2.607'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.9732'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.39%  '
$ws.Range('D37').Value = "'" + '0.5939'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').Value = "'" + '2.680'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('D40').Value = "'" + '1.946'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('D41').Value = "'" + '0.9987'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('D42').Value = "'" + '101.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').Value = "'" + '0.3785'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('D44').Value = "'" + '0.7221'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('D45').Value = "'" + '4.859'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').Value = "'" + '0.05487'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('D47').Value = "'" + '6.239'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.91%  '
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').Value = "'" + '29.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = "'" + '51.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').Value = "'" + '1.000'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.13%  '
